$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 should inherit the same formatting/styles as the row above it (row 11)
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)

# Fill in the new change-log entry
$ws.Range("A12").Value = 41727
$ws.Range("B12").Value = "10"
$ws.Range("C12").Value = "SPA"
$ws.Range("D12").Value = "OS Dispatcher and OSEK functions"
$ws.Range("E12").Value = "In process"

# Move / update the active selection like the authored workbook
$ws.Range("D21").Select()
